# Fix all ICDC breed testcases
# Updates query text for Cases/Samples/Files tabs and adds a new StudyFilesTab row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (CasesTab): update query text (column B) ---
$b2 = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
 WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
WHERE demo.breed IN ['Labrador Retriever']
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS ``Weight (kg)``,
       coalesce(diag.best_response, '') AS ``Response to Treatment``,
       coalesce(co.cohort_description, '') AS ``Cohort``
order by c.case_id asc
limit 100
"@
$ws.Range("B2").Value = $b2

# --- Row 3 (SamplesTab): update query text (column B) ---
$b3 = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Labrador Retriever']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, 
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS ``Sample Site``,
        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,
        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,
        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,
        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,
        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,
        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,
        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``
Order by samp.sample_id LIMIT 100
"@
$ws.Range("B3").Value = $b3

# --- Row 4 (FilesTab): update query text (column B) ---
$b4 = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Labrador Retriever']
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS ``File Name``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_type, '') AS ``File Type``,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(samp.sample_id, '') AS ``Sample ID``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
Order By f.file_name LIMIT 100
"@
$ws.Range("B4").Value = $b4

# --- Row 5 (new StudyFilesTab row) ---
$ws.Range("A5").Value = "StudyFilesTab"
$b5 = @"
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE demo.breed IN ['Labrador Retriever']
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS ``File Name``,
  coalesce(f.file_type, '') AS ``File Type``,
  coalesce("study", '') AS ``Association``,
  coalesce(f.file_description, '') AS ``Description``,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS ``Study Code``
Order By f.file_name LIMIT 100
"@
$ws.Range("B5").Value = $b5
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = $ws.Range("C4").Value2
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value = $ws.Range("D4").Value2
$ws.Range("E5").Value = $ws.Range("E4").Value2

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 288
$ws.Rows.Item(4).RowHeight = 409.6
$ws.Rows.Item(5).RowHeight = 374.4

# --- Column widths: split column B from column C, column B gets its own width ---
$ws.Columns.Item(2).ColumnWidth = 92.44140625

# --- View state: zoom + scroll position + selection ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 55
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C4:E5").Select()
